# Update the "dSF" column (F) values for the specified rows.
# This reflects a "repull data" style update where recalculated values
# replace the previous ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("F7").Value = -1
